$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view changes: scroll/selection moved from column C/N to column F/H ---
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("H28").Select()

# --- Row 20: Navigate to Log in page ---
$ws.Range("P20").Value = "AssertLoginPageUrl"
$ws.Range("R20").Value = "Passed"

# --- Row 21: Navigate to Register page ---
$ws.Range("P21").Value = "AssertRegisterPageUrl"
$ws.Range("R21").Value = "Passed"

# --- Row 22: Navigate to Home page through Logo link ---
$ws.Range("P22").Value = "AssertHomePageUrl"
$ws.Range("R22").Value = "Passed"

# --- Row 23: Log in and enter Account Management ---
$ws.Range("O23").Value = "Change your password"
$ws.Range("P23").Value = "AssertManagePageUrl"
$ws.Range("Q23").Value = "AssertChangePasswordDisplayed"

# --- Row 24: Navigate to create article and click Cancel ---
$ws.Range("P24").Value = "AssertHomePageUrl"

# --- Row 25: Log in and read an article ---
$ws.Range("R25").Value = "Passed"

# --- Row 26: Log in and select article to delete ---
$ws.Range("N26").Value = "http://localhost:60639/Article/Delete/"
$ws.Range("P26").Value = "AssertConfirmDeleteButtonDisplayed"

# --- Row 28: Log in and navigate to password change page ---
$ws.Range("R28").Value = "Passed"
